$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Graduate Students")

$ws.Range("A2").Value = "All students"
$ws.Range("A3").Value = "Male"
$ws.Range("A4").Value = "Female"
$ws.Range("A5").Value = "U.S. citizens and permanent residents"
$ws.Range("A6").Value = "Hispanic or Latino"
$ws.Range("A7").Value = "Not Hispanic or Latino"
$ws.Range("A8").Value = "American Indian or Alaska Native"
$ws.Range("A9").Value = "Asian"
$ws.Range("A10").Value = "Black or African American"
$ws.Range("A11").Value = "Native Hawaiian or Other Pacific Islander"
$ws.Range("A12").Value = "White"
$ws.Range("A13").Value = "More than one race"
$ws.Range("A14").Value = "Unknown ethnicity and race"
$ws.Range("A15").Value = "Temporary visa holders"
$ws.Range("A16").Value = "Science and engineering"
$ws.Range("A17").Value = "Science"
$ws.Range("A18").Value = "Agricultural and veterinary sciences"
$ws.Range("A19").Value = "Biological and biomedical sciences"
$ws.Range("A20").Value = "Communication"
$ws.Range("A21").Value = "Computer and information sciences"
$ws.Range("A22").Value = "Family and consumer sciences and human sciences"
$ws.Range("A23").Value = "Geosciences, atmospheric sciences, and ocean sciences"
$ws.Range("A24").Value = "Mathematics and statistics"
$ws.Range("A25").Value = "Multidisciplinary and interdisciplinary studies"
$ws.Range("A26").Value = "Natural resources and conservation"
$ws.Range("A27").Value = "Psychology"
$ws.Range("A28").Value = "Physical sciences"
$ws.Range("A29").Value = "Social sciences"
$ws.Range("A30").Value = "Engineering"
$ws.Range("A31").Value = "Aerospace, aeronautical, and astronautical engineering"
$ws.Range("A32").Value = "Biological, biomedical, and biosystems engineering"
$ws.Range("A33").Value = "Chemical, petroleum, and chemical-related engineering"
$ws.Range("A34").Value = "Civil, environmental, transportation and related engineering fields"
$ws.Range("A35").Value = "Electrical, electronics, communications and computer engineering"
$ws.Range("A36").Value = "Industrial, manufacturing, systems engineering and operations research"
$ws.Range("A37").Value = "Mechanical engineering"
$ws.Range("A38").Value = "Metallurgical, mining, materials and related engineering fields"
$ws.Range("A39").Value = "Other engineering"
$ws.Range("A40").Value = "Health"
$ws.Range("A41").Value = "Clinical medicine"
$ws.Range("A42").Value = "Other health"
